# #5: cash & deposit done
# Sheet "存款" (deposit): turn row 1 into a real header row and append
# seven new metadata columns (G:M) to every row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("存款")
$landWs = $wb.Worksheets.Item("土地")

# --- Fix row 1: it used to duplicate row 2's data; make it header labels ---
$ws.Range("B1").Value = "bank"
$ws.Range("C1").Value = "deposit_type"
$ws.Range("D1").Value = "currency"
$ws.Range("E1").Value = "owner"
$ws.Range("F1").Value = "total"

# --- New header cells G1:M1, formatted like the existing header (B1) ---
$ws.Range("B1").Copy() | Out-Null
$ws.Range("G1").PasteSpecial(-4122) | Out-Null
$ws.Range("B1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null
$ws.Range("B1").Copy() | Out-Null
$ws.Range("I1").PasteSpecial(-4122) | Out-Null
$ws.Range("B1").Copy() | Out-Null
$ws.Range("J1").PasteSpecial(-4122) | Out-Null
$ws.Range("B1").Copy() | Out-Null
$ws.Range("K1").PasteSpecial(-4122) | Out-Null
$ws.Range("B1").Copy() | Out-Null
$ws.Range("L1").PasteSpecial(-4122) | Out-Null
$ws.Range("B1").Copy() | Out-Null
$ws.Range("M1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("G1").Value = "property_category"
$ws.Range("H1").Value = "category"
$ws.Range("I1").Value = "date"
$ws.Range("J1").Value = "legislator_name"
$ws.Range("K1").Value = "legislator_id"
$ws.Range("L1").Value = "source_file"
$ws.Range("M1").Value = "index"

# Row 2..11 source data that already lives in B:F, used to build the new
# metadata columns G:M for each row.
$rowIndex = @(47, 48, 49, 50, 51, 52, 53, 54, 55, 56)

for ($r = 2; $r -le 11; $r++) {
    # Copy formatting for the new cells from the existing data cell in the
    # same row (B<r>) so the new columns match the rest of the row.
    $ws.Range("B$r").Copy() | Out-Null
    $ws.Range("G$r").PasteSpecial(-4122) | Out-Null
    $ws.Range("B$r").Copy() | Out-Null
    $ws.Range("H$r").PasteSpecial(-4122) | Out-Null
    $ws.Range("B$r").Copy() | Out-Null
    $ws.Range("I$r").PasteSpecial(-4122) | Out-Null
    $ws.Range("B$r").Copy() | Out-Null
    $ws.Range("J$r").PasteSpecial(-4122) | Out-Null
    $ws.Range("B$r").Copy() | Out-Null
    $ws.Range("K$r").PasteSpecial(-4122) | Out-Null
    $ws.Range("B$r").Copy() | Out-Null
    $ws.Range("L$r").PasteSpecial(-4122) | Out-Null
    $ws.Range("B$r").Copy() | Out-Null
    $ws.Range("M$r").PasteSpecial(-4122) | Out-Null
    $excel.CutCopyMode = $false

    $ws.Range("G$r").Value = "deposit"
    $ws.Range("H$r").Value = "normal"

    # "2011-11-21" must stay a text string, not become a date serial. Pull
    # it as a value (not a formula) from 土地!K2, which already stores the
    # identical text "2011-11-21" as a string, so no number-format change
    # (and therefore no new style) is needed on the destination cell.
    $landWs.Range("K2").Copy() | Out-Null
    $ws.Range("I$r").PasteSpecial(-4163) | Out-Null
    $excel.CutCopyMode = $false

    $ws.Range("J$r").Value = "盧秀燕"
    $ws.Range("K$r").Value = 869
    $ws.Range("L$r").Value = "tmp9eb41"
    $ws.Range("M$r").Value = $rowIndex[$r - 2]
}
